$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44587
$ws.Range("M2").Value = 165
$ws.Range("N2").Value = 6500
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6742
$ws.Range("S2").Value = 3371

# Row 3
$ws.Range("D3").Value = 44586
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 3500

# Row 4
$ws.Range("D4").Value = 44588
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6750
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 3375

# Row 5
$ws.Range("D5").Value = 44627
$ws.Range("M5").Value = 45

# Row 6
$ws.Range("D6").Value = 44606
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("R6").Value = "Provincia de Linares"
$ws.Range("S6").Value = 3500

# Row 7
$ws.Range("D7").Value = 44628
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("R7").Value = "Provincia de Linares"
$ws.Range("S7").Value = 3000

# Row 9
$ws.Range("D9").Value = 44214
$ws.Range("M9").Value = 48

# Row 10
$ws.Range("D10").Value = 44211
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("S10").Value = 3000

# Row 11
$ws.Range("D11").Value = 44589
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 6000
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 3000

# Row 13
$ws.Range("D13").Value = 44614
$ws.Range("R13").Value = "Provincia de Linares"

# Row 14
$ws.Range("D14").Value = 44592
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value = 4000

# Row 15
$ws.Range("D15").Value = 44582
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6500
$ws.Range("P15").Value = 6233
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 3116
